$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44657
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("S2").Value = 1857

# Row 3
$ws.Range("D3").Value = 44321
$ws.Range("M3").Value = 140
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11500
$ws.Range("S3").Value = 1643

# Row 4
$ws.Range("D4").Value = 44321
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("S4").Value = 1143

# Row 5
$ws.Range("D5").Value = 44315
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("S5").Value = 2000

# Row 6
$ws.Range("D6").Value = 44315
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("S6").Value = 1786

# Row 7
$ws.Range("D7").Value = 44315
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 10000
$ws.Range("P7").Value = 10500
$ws.Range("S7").Value = 1500

# Row 8
$ws.Range("D8").Value = 44344
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 1714

# Row 9
$ws.Range("D9").Value = 44306
$ws.Range("L9").Value = 'Primera'
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10
$ws.Range("D10").Value = 44306
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 9000
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 1286

# Row 11
$ws.Range("D11").Value = 44322
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 11000
$ws.Range("S11").Value = 1571

# Row 12
$ws.Range("D12").Value = 44314
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("S12").Value = 1857

# Row 13
$ws.Range("D13").Value = 44314
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 11000
$ws.Range("O13").Value = 11000
$ws.Range("P13").Value = 11000
$ws.Range("S13").Value = 1571

# Row 14
$ws.Range("D14").Value = 44302
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 340
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 12500
$ws.Range("R14").Value = 'Provincia de Santiago'
$ws.Range("S14").Value = 1786

# Row 15
$ws.Range("D15").Value = 44316
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("S15").Value = 1857

# Row 16
$ws.Range("D16").Value = 44316
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 11000
$ws.Range("S16").Value = 1571

# Row 17
$ws.Range("D17").Value = 44623
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 2286

# Row 18
$ws.Range("D18").Value = 44643
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 2143

# Row 19
$ws.Range("D19").Value = 44300
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 12000
$ws.Range("P19").Value = 12500
$ws.Range("R19").Value = 'Provincia de Santiago'
$ws.Range("S19").Value = 1786

# Row 20
$ws.Range("D20").Value = 44644
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 85
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("S20").Value = 2000

# Row 21
$ws.Range("D21").Value = 44312
$ws.Range("L21").Value = 'Primera'
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("S21").Value = 1857

# Row 22
$ws.Range("D22").Value = 44312
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 11000
$ws.Range("O22").Value = 11000
$ws.Range("P22").Value = 11000
$ws.Range("S22").Value = 1571

# Row 23
$ws.Range("D23").Value = 44342
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("S23").Value = 1714

# Row 25
$ws.Range("D25").Value = 44335
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("S25").Value = 2000

# Row 26
$ws.Range("D26").Value = 44349
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 70
$ws.Range("N26").Value = 12000
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 12000
$ws.Range("S26").Value = 1714

# Row 27
$ws.Range("D27").Value = 44307
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 70
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("S27").Value = 2000

# New row 28 (newly recorded entry)
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C28").Value = 'Metropolitana'
$ws.Range("D28").Value = 44307
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 'Fruta'
$ws.Range("G28").Value = 100101
$ws.Range("H28").Value = 'Berries'
$ws.Range("I28").Value = 100101006
$ws.Range("J28").Value = 'Higo'
$ws.Range("K28").Value = 'Sin especificar'
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("Q28").Value = '$/bandeja 7 kilos'
$ws.Range("R28").Value = 'Región Metropolitana'
$ws.Range("S28").Value = 1429
$ws.Range("T28").Value = 7

# New row inherits the same date number format used throughout column D
$ws.Range("D28").NumberFormat = $ws.Range("D27").NumberFormat
